# Actualizar 02-05-2021 19-57-53
# Appends a fresh block of 14 availability-check rows (772-785) to Sheet1,
# mirroring the existing 758-771 block with a newer check timestamp, and
# nudges the previous block's timestamp (44232.81028977159 -> 44232.81028976852).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Correct the timestamp already stored for the 758-771 block ---------
$prevTimestamp = 44232.81028976852
for ($r = 758; $r -le 771; $r++) {
    $ws.Cells.Item($r, 4).Value = $prevTimestamp
}

# --- 2. Append the new block (rows 772-785) ---------------------------------
$newTimestamp = 44232.83168540164

$newRows = @(
    @{ Row = 772; Name = "Odoo";               Url = "https://www.dataintelligence-group.com/";                      SubAddress = $null },
    @{ Row = 773; Name = "Blackbox";            Url = "https://serviciodashboard.azurewebsites.net/";                 SubAddress = $null },
    @{ Row = 774; Name = "PowerBI";             Url = "https://powerbi.microsoft.com/es-es/";                        SubAddress = $null },
    @{ Row = 775; Name = "Dropbox";             Url = "https://www.dropbox.com/";                                    SubAddress = $null },
    @{ Row = 776; Name = "Odoo";                Url = "https://dataintelligence.store/";                             SubAddress = $null },
    @{ Row = 777; Name = "GEE";                 Url = "https://app-data-i.users.earthengine.app/";                   SubAddress = $null },
    @{ Row = 778; Name = "UtilidadesOdoo";      Url = "https://odooutil.azurewebsites.net/";                         SubAddress = $null },
    @{ Row = 779; Name = "Filtros Dashboard";   Url = "https://filtradordashboard.azurewebsites.net/";                SubAddress = $null },
    @{ Row = 780; Name = "MapStore";            Url = "https://ide.dataintelligence-group.com/mapstore/";            SubAddress = "/"; Display = "https://ide.dataintelligence-group.com/mapstore/#/" },
    @{ Row = 781; Name = "GeoServer";           Url = "https://ide.dataintelligence-group.com/geoserver/web/?0";      SubAddress = $null },
    @{ Row = 782; Name = "Tomcat";              Url = "https://ide.dataintelligence-group.com/";                     SubAddress = $null },
    @{ Row = 783; Name = "Shiny";               Url = "https://rpubs.com/dataintelligence/";                         SubAddress = $null },
    @{ Row = 784; Name = "Github";              Url = "https://github.com/Sud-Austral/";                             SubAddress = $null },
    @{ Row = 785; Name = "EZ Exporter";         Url = "https://ezexporter.highviewapps.com/exports/export-profile/"; SubAddress = $null }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value = $entry.Name
    if ($entry.ContainsKey("Display")) {
        $ws.Cells.Item($r, 2).Value = $entry.Display
    } else {
        $ws.Cells.Item($r, 2).Value = $entry.Url
    }
    $ws.Cells.Item($r, 3).Value = "Disponible"
    $ws.Cells.Item($r, 4).Value = $newTimestamp

    $linkCell = $ws.Cells.Item($r, 2)
    if ($entry.SubAddress) {
        $ws.Hyperlinks.Add($linkCell, $entry.Url, $entry.SubAddress)
    } else {
        $ws.Hyperlinks.Add($linkCell, $entry.Url)
    }
    # Hyperlinks.Add() re-derives the cell style; pin it back to the shared
    # "Hyperlink" cell style so it reuses the same style index as every
    # other link cell instead of allocating a near-duplicate one.
    $linkCell.Style = "Hyperlink"

    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
